$d = $word.ActiveDocument

$d.Content.Find.Execute("vamos a crear un ejemplo de", $true, $false, $false, $false, $false, $true, 1, $false, "se construirá un aplicación web de ejemplo de", 2)
